# Update "Lương" worksheet with new rows / values for Luong ca nhan calculation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# --- Update existing summary values ---
$ws.Range("B2").Value = 14
$ws.Range("B3").Value = 490000

# --- Insert "Ứng lương tại CẦN THƠ" row after current row 10 (before LONG XUYÊN block) ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$ws.Range("B11").Value = 0

# --- Insert "Ứng lương tại LONG XUYÊN" row after LONG XUYÊN's "Công phụ phẫu 2" (now row 18) ---
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$ws.Range("B19").Value = 0

# --- Update SÓC TRĂNG "Lương cơ bản" value (now row 20) ---
$ws.Range("B20").Value = 1500000

# --- Insert "Ứng lương tại SÓC TRĂNG" row after SÓC TRĂNG's "Công phụ phẫu 2" (now row 26) ---
$ws.Rows.Item(27).Insert()
$ws.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$ws.Range("B27").Value = 0

# --- Append the four new "Tổng lương" summary rows ---
$ws.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$ws.Range("B29").Value = 0

$ws.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$ws.Range("B30").Value = 5665000

$ws.Range("A31").Value = "Tổng lương"
$ws.Range("B31").Value = 5665000
